$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: target row -> source row (permutation of the weekly data rows 2-61)
$mapping = @{
    2 = 39
    3 = 35
    4 = 51
    5 = 48
    6 = 7
    7 = 19
    8 = 4
    9 = 11
    10 = 32
    11 = 38
    12 = 45
    13 = 59
    14 = 46
    15 = 40
    16 = 24
    17 = 57
    18 = 2
    19 = 34
    20 = 26
    21 = 49
    22 = 37
    23 = 13
    24 = 41
    25 = 10
    26 = 52
    27 = 56
    28 = 16
    29 = 18
    30 = 33
    31 = 6
    32 = 28
    33 = 61
    34 = 5
    35 = 21
    36 = 44
    37 = 22
    38 = 43
    39 = 30
    40 = 12
    41 = 25
    42 = 29
    43 = 23
    44 = 15
    45 = 60
    46 = 50
    47 = 20
    48 = 54
    49 = 31
    50 = 9
    51 = 8
    52 = 3
    53 = 36
    54 = 27
    55 = 17
    56 = 14
    57 = 42
    58 = 47
    59 = 55
    60 = 53
    61 = 58
}

$cols = @("D", "J", "K", "L", "M", "P")

# Step 1: snapshot all current (before) values for the affected columns/rows
$before = @{}
for ($r = 2; $r -le 61; $r++) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $before[$r] = $rowVals
}

# Step 2: write back the permuted values
for ($t = 2; $t -le 61; $t++) {
    $s = $mapping[$t]
    $srcVals = $before[$s]
    foreach ($col in $cols) {
        $ws.Range("$col$t").Value2 = $srcVals[$col]
    }
}

Write-Host "Permutation applied."